$wb = $excel.ActiveWorkbook

# 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1386
$ws.Range("F4").Value = 2107
$ws.Range("F5").Value = 5734
$ws.Range("F9").Value = 6562
$ws.Range("F10").Value = 196
$ws.Range("F22").Value = 942
$ws.Range("F23").Value = 316
$ws.Range("F32").Value = 18
$ws.Range("F33").Value = 271

# 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F21").Value = 130

# 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 536
$ws.Range("F8").Value = 797

# 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1386
$ws.Range("F8").Value = 536
$ws.Range("F9").Value = 536
$ws.Range("F14").Value = 5734
$ws.Range("F19").Value = 6562
$ws.Range("F20").Value = 196
$ws.Range("F28").Value = 130
$ws.Range("F32").Value = 942
$ws.Range("F40").Value = 18
$ws.Range("F42").Value = 271
